$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 2 Correspond Handoff/Handback Datetime for the
# f6168691-...15fcf0ac...zh-cn.xlf handback entry
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 12:55:03"
$wsZhCn.Range("H2").Value = "2016-03-20 12:55:22"

# de-de sheet: row 2 Correspond Handoff/Handback Datetime for the
# f6168691-...15fcf0ac...de-de.xlf handback entry
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 12:55:07"
$wsDeDe.Range("H2").Value = "2016-03-20 12:55:27"
